# "change zorder of lines in plot"
# The "zorder" column (S) for every data row (2-36) is bumped from 1 to 2,
# and the sheet's view/selection is moved (pane scrolled to K2, active
# selection moved to T10).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# --- bump the zorder column (S2:S36) from 1 to 2 ---------------------------
for ($r = 2; $r -le 36; $r++) {
    $ws.Cells.Item($r, 19).Value = 2
}

# --- update the sheet view: scroll position + active selection ------------
$ws.Activate()

$win = $excel.ActiveWindow
$win.ScrollColumn = 11   # column K
$win.ScrollRow = 2

$ws.Range("T10").Select() | Out-Null
